$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by copying the structurally-identical
#    "2022-Q1" sheet (same column layout/styles), placed right after the
#    "总计" summary sheet, then overwrite its data.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$q1Sheet.Copy($null, $totalSheet)
$newSheet = $wb.Worksheets.Item("2022-Q1 (2)")
$newSheet.Name = "2022-Q4"

# Clear out the copied data rows (2..32) so no stale 2022-Q1 values remain.
$newSheet.Range("A2:H32").ClearContents()

$fundData = @(
    ,@(0, '159941', '广发纳斯达克100ETF（QDII）', '114.77', '90.42', '11.26', '12.9231', 1)
    ,@(1, '513100', '国泰纳斯达克100（QDII-ETF）', '51.50', '90.72', '11.68', '6.0152', 1)
    ,@(2, '513500', '博时标普500ETF（QDII）', '93.59', '95.36', '5.77', '5.4001', 1)
    ,@(3, '159632', '华安纳斯达克100ETF（QDII）', '30.39', '87.77', '10.34', '3.1423', 2)
    ,@(4, '160213', '国泰纳斯达克100指数（QDII）', '15.65', '90.80', '10.67', '1.6699', 2)
    ,@(5, '000834', '大成纳斯达克100指数（QDII）', '15.61', '81.77', '9.59', '1.4970', 2)
    ,@(6, '513300', '华夏纳斯达克100ETF（QDII）', '12.43', '97.54', '11.47', '1.4257', 1)
    ,@(7, '000043', '嘉实美国成长股票（QDII）人民币', '12.69', '92.23', '8.85', '1.1231', 1)
    ,@(8, '000044', '嘉实美国成长股票（QDII）美元现汇', '12.69', '92.23', '8.85', '1.1231', 1)
    ,@(9, '012868', '易方达标普信息科技指数（QDII-LOF）人民币 C', '5.09', '91.36', '21.46', '1.0923', 1)
    ,@(10, '161128', '易方达标普信息科技指数（QDII-LOF）人民币', '5.09', '91.36', '21.46', '1.0923', 1)
    ,@(11, '003721', '易方达标普信息科技指数（QDII-LOF）美元A', '4.93', '91.36', '21.46', '1.0580', 1)
    ,@(12, '161130', '易方达纳斯达克100指数人民币（QDII-LOF）', '7.77', '90.34', '10.62', '0.8252', 2)
    ,@(13, '003722', '易方达纳斯达克100指数美元（QDII-LOF）A', '7.77', '90.34', '10.62', '0.8252', 2)
    ,@(14, '000041', '华夏全球精选股票（QDII）', '18.44', '85.51', '2.57', '0.4739', 6)
    ,@(15, '015203', '汇添富全球移动互联灵活配置混合（QDII）D', '11.52', '92.14', '2.56', '0.2949', 9)
    ,@(16, '001668', '汇添富全球移动互联灵活配置混合（QDII）A', '11.48', '92.14', '2.56', '0.2939', 9)
    ,@(17, '012860', '易方达标普500指数（QDII-LOF）人民币 C', '4.75', '91.65', '5.52', '0.2622', 1)
    ,@(18, '161125', '易方达标普500指数（QDII-LOF）人民币', '4.75', '91.65', '5.52', '0.2622', 1)
    ,@(19, '003718', '易方达标普500指数（QDII-LOF）美元A', '4.65', '91.65', '5.52', '0.2567', 1)
    ,@(20, '016532', '嘉实纳斯达克100指数（QDII）A人民币', '1.12', '94.67', '11.14', '0.1248', 2)
    ,@(21, '016533', '嘉实纳斯达克100指数（QDII）C人民币', '1.12', '94.67', '11.14', '0.1248', 2)
    ,@(22, '016534', '嘉实纳斯达克100指数（QDII）A美元现汇', '1.12', '94.67', '11.14', '0.1248', 2)
    ,@(23, '016535', '嘉实纳斯达克100指数（QDII）C美元现汇', '1.12', '94.67', '11.14', '0.1248', 2)
    ,@(24, '016055', '博时纳斯达克100指数（QDII）A人民币', '1.06', '90.62', '10.65', '0.1129', 2)
    ,@(25, '016057', '博时纳斯达克100指数（QDII）C人民币', '1.06', '90.62', '10.65', '0.1129', 2)
    ,@(26, '016056', '博时纳斯达克100指数（QDII）A美元现汇', '1.06', '90.62', '10.65', '0.1129', 2)
    ,@(27, '016058', '博时纳斯达克100指数（QDII）C美元现汇', '1.06', '90.62', '10.65', '0.1129', 2)
    ,@(28, '013328', '嘉实全球价值股票（QDII）人民币', '1.62', '90.66', '3.96', '0.0642', 1)
    ,@(29, '013329', '嘉实全球价值股票（QDII）美元现汇', '1.62', '90.66', '3.96', '0.0642', 1)
    ,@(30, '159612', '国泰标普500ETF（QDII）', '0.86', '94.21', '5.85', '0.0503', 1)
    ,@(31, '005698', '华夏全球科技先锋混合（QDII）', '0.60', '83.35', '8.18', '0.0491', 2)
    ,@(32, '012869', '易方达标普信息科技指数（QDII-LOF）美元 C', '0.16', '91.36', '21.46', '0.0343', 1)
    ,@(33, '012870', '易方达纳斯达克100指数人民币（QDII-LOF）C', '0.21', '90.34', '10.62', '0.0223', 2)
    ,@(34, '012871', '易方达纳斯达克100指数美元（QDII-LOF）C', '0.21', '90.34', '10.62', '0.0223', 2)
    ,@(35, '015205', '银华全球新能源车量化优选股票（QDII）C', '0.21', '86.58', '5.97', '0.0125', 3)
    ,@(36, '159655', '华夏标普500ETF（QDII）', '0.21', '93.70', '5.67', '0.0119', 1)
    ,@(37, '015204', '银华全球新能源车量化优选股票（QDII）A', '0.16', '86.58', '5.97', '0.0096', 3)
    ,@(38, '014002', '浦银安盛全球智能科技股票（QDII）C', '0.30', '42.55', '2.35', '0.0070', 2)
    ,@(39, '006555', '浦银安盛全球智能科技股票（QDII）A', '0.25', '42.55', '2.35', '0.0059', 2)
    ,@(40, '012861', '易方达标普500指数（QDII-LOF）美元 C', '0.10', '91.65', '5.52', '0.0055', 1)
    ,@(41, '015202', '汇添富全球移动互联灵活配置混合（QDII）C', '0.01', '92.14', '2.56', '0.0003', 9)
)


# Make sure column A (the numeric row index) carries the same style as the
# rest of that column for the newly-needed rows (33..43) by extending the
# existing formatted cell (A32) downward before writing values.
$newSheet.Range("A32").Copy()
$newSheet.Range("A33:A43").PasteSpecial(-4122)

$r = 2
foreach ($row in $fundData) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]

    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 2).Style = "Normal"

    $newSheet.Cells.Item($r, 3).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 3).Style = "Normal"

    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 4).Style = "Normal"

    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 5).Style = "Normal"

    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 6).Style = "Normal"

    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 7).Style = "Normal"

    $newSheet.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q4 at
#    the top of the data (row 2), pushing the rest down, then renumber
#    the index column (A) sequentially.
# ---------------------------------------------------------------------
$ws = $totalSheet
$ws.Rows.Item(2).Insert()
$ws.Range("A3:D3").Copy()
$ws.Range("A2:D2").PasteSpecial(-4122)

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "2022-Q4"
$ws.Cells.Item(2, 3).Value = 42
$ws.Cells.Item(2, 4).Value = 42.37

# Renumber the index column for the rows that shifted down.
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(7, 1).Value = 5

Write-Host "Done"
